$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "Ord12992018"
$ws.Range("O2").Value = "Ser12992018"
$ws.Range("P2").Value = "Test1234"
$ws.Range("R2").Value = "Repair And Replace"
$ws.Range("S2").Value = "Received"

# --- Row 5 ---
$ws.Range("A5").Value = "SN00006"
$ws.Range("C5").Value = "ZBSN1880"
$ws.Range("F5").Value = "TestPrinter"
$ws.Range("G5").Value = "TestBin"
$ws.Range("H5").Value = "TestBin"

# --- Row 7 ---
$ws.Range("A7").Value = "SN00009"
$ws.Range("B7").Value = "TestBin"
$ws.Range("D7").Value = "Battery Status Test "
$ws.Range("E7").Value = "TestBin"

# --- Row 9 ---
$ws.Range("A9").Value = "AH00001"
$ws.Range("B9").Value = "TestBin"
$ws.Range("C9").Value = "Repairable"
$ws.Range("D9").Value = "test123"
$ws.Range("E9").Value = "Battery-Diminished"
$ws.Range("F9").Value = "Battery Tray"
$ws.Range("G9").Value = "Test Repairable"
$ws.Range("H9").Value = "Test NTF"
$ws.Range("I9").Value = "Test Unrepairable"
$ws.Range("J9").Value = "REPLACE/INSTALL"
$ws.Range("K9").Value = "Battery-Diminished"

# --- Row 11 ---
$ws.Range("A11").Value = "AH00001"
$ws.Range("B11").Value = "TestBin"
$ws.Range("C11").Value = "PASS"
$ws.Range("D11").Value = "Aesthetic Check"

# --- Row 13 ---
$ws.Range("A13").Value = "Pack878899"
$ws.Range("B13").Value = "AG00001"
$ws.Range("C13").Value = "AH00001"

# --- Row 15 ---
$ws.Range("A15").Value = "Pack878899"
$ws.Range("B15").Value = "AH00001"
$ws.Range("C15").Value = "External"
$ws.Range("D15").Value = "Test2134"
$ws.Range("E15").Value = "FEDEX"
$ws.Range("F15").Value = "STANDARD OVERNIGHT"
$ws.Range("G15").Value = 1000
$ws.Range("H15").Value = "USD"

# --- Row 16 ---
$ws.Range("A16").Value = "Create Part"

# --- Row 17 ---
$ws.Range("A17").Value = "Ravi1234"
$ws.Range("B17").Value = "TestDescription"
$ws.Range("F17").Value = "test123"

# --- Row 18 ---
$ws.Range("A18").Value = "Component Picking"

# --- Row 19 ---
$ws.Range("A19").Value = "C01A01 "
$ws.Range("B19").Value = "test123"
$ws.Range("C19").Value = "Ser12992018"

# --- Row 20 ---
$ws.Range("A20").Value = "Component Putaway"

# --- Row 21 ---
$ws.Range("A21").Value = "test123"
$ws.Range("B21").Value = "Ser12992018"

# --- Row 22 ---
$ws.Range("A22").Value = "Part Picking"

# --- Row 23 ---
$ws.Range("A23").Value = "test123"
$ws.Range("B23").Value = "C01A01 "
$ws.Range("C23").Value = "PartsPicking "

# --- Row 24 ---
$ws.Range("A24").Value = "Deliver Part"

# --- Row 25 ---
$ws.Range("B25").Value = "Ser12992018"

# --- Row 26 ---
$ws.Range("A26").Value = "DueDateCalc"

# --- Row 27 ---
$ws.Range("B27").Value = "HWPN1234"
$ws.Range("D27").Value = "2N10141077"

# --- Row 28 ---
$ws.Range("A28").Value = "Manage Region"

# --- Row 29 ---
$ws.Range("A29").Value = "Bengal"
$ws.Range("B29").Value = "India Standard Time"
$ws.Range("C29").Value = "Dateline Standard Time"

# --- Row 30 ---
$ws.Range("A30").Value = "BranchMaster"

# --- Row 31 ---
$ws.Range("A31").Value = "New York"
$ws.Range("B31").Value = "TestBranch"
$ws.Range("C31").Value = "TestAdd"
$ws.Range("D31").Value = "Address line 2"
$ws.Range("E31").Value = "TestCity"
$ws.Range("F31").Value = "TestState"
$ws.Range("G31").Value = "TestZip"
$ws.Range("H31").Value = "TestCountry"
$ws.Range("I31").Value = "India Standard Time"
$ws.Range("J31").Value = "Contact Name"
$ws.Range("L31").Value = "Test@ctdi.com"

# --- Row 32 ---
$ws.Range("A32").Value = "USerMaster"

# --- Row 33 ---
$ws.Range("A33").Value = "Test"
$ws.Range("B33").Value = "User"
$ws.Range("C33").Value = "Tuser"
$ws.Range("D33").Value = "Tuser@gmail.com"
$ws.Range("E33").Value = "Developer"
$ws.Range("F33").Value = "Test"
$ws.Range("G33").Value = "ViewOrders"

# --- Row 34 ---
$ws.Range("A34").Value = "BinsMaster"

# --- Row 35 ---
$ws.Range("A35").Value = "Test"
$ws.Range("B35").Value = "TestDescription"
$ws.Range("C35").Value = "Standard"
$ws.Range("E35").Value = "Aisle1373A"
$ws.Range("G35").Value = "EditDescription"

# --- Column widths: merge F & G into a single 22.28515625-wide block ---
$ws.Columns("F:G").ColumnWidth = 22.28515625

# --- View: scroll back to top and select H15 ---
$ws.Range("A1").Select()
$ws.Range("H15").Select()
